$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 3323.7273
$ws.Range("I113").Value = 2456.3333
$ws.Range("J113").Value = 3649
$ws.Range("K113").Value = 2456.3333
$ws.Range("L113").Value = 3649
$ws.Range("M113").Value = 797.6667000000002
$ws.Range("N113").Value = -10157
$ws.Range("H116").Value = 2073.4707
$ws.Range("I116").Value = 2000
$ws.Range("J116").Value = 2156.125
$ws.Range("K116").Value = 2000
$ws.Range("L116").Value = 2156.125
$ws.Range("M116").Value = 1442
$ws.Range("N116").Value = -9040.125
$ws.Range("H138").Value = 6004.907
$ws.Range("I138").Value = 1585.4
$ws.Range("J138").Value = 8604.617
$ws.Range("K138").Value = 4756.200000000001
$ws.Range("L138").Value = 25813.851
$ws.Range("M138").Value = 383.7999999999993
$ws.Range("N138").Value = -36093.851

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 835.8823
$ws.Range("I97").Value = 763.125
$ws.Range("K97").Value = 763.125
$ws.Range("M97").Value = -267.125
$ws.Range("H102").Value = 3785.7144
$ws.Range("I102").Value = 3618.182
$ws.Range("K102").Value = 3618.182
$ws.Range("M102").Value = -1996.182
$ws.Range("H122").Value = 2242.4285
$ws.Range("I122").Value = 2114.2856
$ws.Range("J122").Value = 2370.5715
$ws.Range("K122").Value = 6342.8568
$ws.Range("L122").Value = 7111.7145
$ws.Range("M122").Value = -3892.8568
$ws.Range("N122").Value = -12011.7145

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 836.6070999999999
$ws.Range("J94").Value = 1263.75
$ws.Range("L94").Value = 1263.75
$ws.Range("N94").Value = -2165.75
$ws.Range("H105").Value = 4928.7144
$ws.Range("I105").Value = 3163.3333
$ws.Range("K105").Value = 3163.3333
$ws.Range("M105").Value = -1416.3333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 35901.57
$ws.Range("J82").Value = 35901.57
$ws.Range("L82").Value = 35901.57
$ws.Range("N82").Value = -36667.57
$ws.Range("H85").Value = 35901.57
$ws.Range("J85").Value = 35901.57
$ws.Range("L85").Value = 35901.57
$ws.Range("N85").Value = -38553.57
$ws.Range("H97").Value = 885.39026
$ws.Range("I97").Value = 844.37036
$ws.Range("J97").Value = 964.5
$ws.Range("K97").Value = 844.37036
$ws.Range("L97").Value = 964.5
$ws.Range("M97").Value = -348.37036
$ws.Range("N97").Value = -1956.5
$ws.Range("H126").Value = 2214.6333
$ws.Range("I126").Value = 1738.8948
$ws.Range("J126").Value = 3036.3635
$ws.Range("K126").Value = 5216.6844
$ws.Range("L126").Value = 9109.0905
$ws.Range("M126").Value = -2746.6844
$ws.Range("N126").Value = -14049.0905

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5388.3125
$ws.Range("I7").Value = 5400.8
$ws.Range("J7").Value = 5367.5
$ws.Range("K7").Value = 5400.8
$ws.Range("L7").Value = 5367.5
$ws.Range("M7").Value = -5288.8
$ws.Range("N7").Value = -5591.5
$ws.Range("H68").Value = 4585.857
$ws.Range("I68").Value = 4400.2
$ws.Range("J68").Value = 5050
$ws.Range("K68").Value = 4400.2
$ws.Range("L68").Value = 5050
$ws.Range("M68").Value = -3651.2
$ws.Range("N68").Value = -6548
$ws.Range("H71").Value = 4585.857
$ws.Range("I71").Value = 4400.2
$ws.Range("J71").Value = 5050
$ws.Range("K71").Value = 22001
$ws.Range("L71").Value = 25250
$ws.Range("M71").Value = -18257
$ws.Range("N71").Value = -32738
$ws.Range("H93").Value = 2053
$ws.Range("J93").Value = 2328.3333
$ws.Range("L93").Value = 2328.3333
$ws.Range("N93").Value = -4824.3333
$ws.Range("H94").Value = 12296
$ws.Range("J94").Value = 12296
$ws.Range("L94").Value = 12296
$ws.Range("N94").Value = -13648
$ws.Range("H96").Value = 67000
$ws.Range("J96").Value = 67000
$ws.Range("L96").Value = 67000
$ws.Range("N96").Value = -72492
$ws.Range("H99").Value = 36500
$ws.Range("J99").Value = 53000
$ws.Range("L99").Value = 53000
$ws.Range("N99").Value = -58990
$ws.Range("H103").Value = 50000
$ws.Range("J103").Value = 50000
$ws.Range("L103").Value = 50000
$ws.Range("N103").Value = -52344
$ws.Range("H104").Value = 33456.668
$ws.Range("J104").Value = 33456.668
$ws.Range("L104").Value = 33456.668
$ws.Range("N104").Value = -40444.668
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = $null
$ws.Range("H122").Value = 7101.4546
$ws.Range("I122").Value = 6673.1113
$ws.Range("K122").Value = 20019.3339
$ws.Range("M122").Value = -17569.3339
$ws.Range("H126").Value = 5388.3125
$ws.Range("I126").Value = 5400.8
$ws.Range("J126").Value = 5367.5
$ws.Range("K126").Value = 16202.4
$ws.Range("L126").Value = 16102.5
$ws.Range("M126").Value = -13732.4
$ws.Range("N126").Value = -21042.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H92").Value = 29000
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 29000
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 29000
$ws.Range("M92").Value = $null
$ws.Range("N92").Value = -33992
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = $null
$ws.Range("H95").Value = 42500
$ws.Range("J95").Value = 42500
$ws.Range("L95").Value = 42500
$ws.Range("N95").Value = -47992
$ws.Range("H97").Value = 10000
$ws.Range("J97").Value = 10000
$ws.Range("L97").Value = 10000
$ws.Range("N97").Value = -11982
$ws.Range("H98").Value = 24266.334
$ws.Range("J98").Value = 24266.334
$ws.Range("L98").Value = 24266.334
$ws.Range("N98").Value = -30256.334
$ws.Range("H99").Value = 65000
$ws.Range("J99").Value = 65000
$ws.Range("L99").Value = 65000
$ws.Range("N99").Value = -70990
$ws.Range("H100").Value = 3075.75
$ws.Range("I100").Value = 800
$ws.Range("K100").Value = 1600
$ws.Range("M100").Value = -1059
$ws.Range("H101").Value = 26339.572
$ws.Range("J101").Value = 26339.572
$ws.Range("L101").Value = 26339.572
$ws.Range("N101").Value = -32829.572
$ws.Range("H102").Value = 59000
$ws.Range("J102").Value = 59000
$ws.Range("L102").Value = 59000
$ws.Range("N102").Value = -65490
$ws.Range("H103").Value = 46801.332
$ws.Range("J103").Value = 46801.332
$ws.Range("L103").Value = 46801.332
$ws.Range("N103").Value = -49145.332
$ws.Range("H104").Value = 0
$ws.Range("J104").Value = 0
$ws.Range("L104").Value = 0
$ws.Range("N104").Value = $null
$ws.Range("H105").Value = 39807.5
$ws.Range("J105").Value = 39807.5
$ws.Range("L105").Value = 39807.5
$ws.Range("N105").Value = -46795.5
$ws.Range("H106").Value = 55000
$ws.Range("J106").Value = 55000
$ws.Range("L106").Value = 55000
$ws.Range("N106").Value = -57524
$ws.Range("H122").Value = 28002.5
$ws.Range("J122").Value = 28002.5
$ws.Range("L122").Value = 84007.5
$ws.Range("N122").Value = -88907.5
$ws.Range("H126").Value = 1171
$ws.Range("I126").Value = 1267.5264
$ws.Range("J126").Value = 967.2222
$ws.Range("K126").Value = 3802.5792
$ws.Range("L126").Value = 2901.6666
$ws.Range("M126").Value = -1332.5792
$ws.Range("N126").Value = -7841.6666
$ws.Range("H133").Value = 40810
$ws.Range("J133").Value = 40810
$ws.Range("L133").Value = 40810
$ws.Range("N133").Value = -50930
